$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: expand the single "Food Safety Course" row into three levels
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Rename the existing "Food Safety Course" entry to "Level 1"
$ws1.Range("A2").Value = "Food Safety Course Level 1"

# Insert two new rows right after A2 (before the old "BCSS Course" row)
# so that BCSS Course is pushed down from row 3 to row 5.
$ws1.Range("A3").EntireRow.Insert()
$ws1.Range("A3").EntireRow.Insert()

$ws1.Range("A3").Value = "Food Safety Course Level 2"
$ws1.Range("A4").Value = "Food Safety Course Level 3"
# A5 already holds "BCSS Course" after the shift, nothing to change there.

# Update the selected cell shown in the sheet view
$ws1.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------
# AKC Rankings: replace the 10 stale "Food Safety Course" rows + the
# trailing "BCSS Course" row with 4 fresh rows (3 levels + BCSS Course)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("AKC Rankings")

# Drop the old data rows 6-11, keeping rows 1 (header) through 5.
$ws2.Range("A6:A11").EntireRow.Delete() | Out-Null

function Set-RankingRow($row, $term, $rank, $date) {
    $ws2.Cells.Item($row, 1).Value = $term

    $ws2.Cells.Item($row, 2).Value = $rank

    $dateCell = $ws2.Cells.Item($row, 3)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $date
    $dateCell.Style = "Normal"
}

Set-RankingRow 2 "Food Safety Course Level 1" 10 "2025-11-06"
Set-RankingRow 3 "Food Safety Course Level 2" 7  "2025-11-06"
Set-RankingRow 4 "Food Safety Course Level 3" 14 "2025-11-06"
Set-RankingRow 5 "BCSS Course" 2 "2025-11-06"
